$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S")

$data = @"
n:7|s:model_24_9_8|n:0.9978655741283762|n:0.9947922486835233|n:0.9999999999638428|n:0.9999999999982546|n:0.9999999999956094|n:805.7457665380091|n:1965.92612196995|n:6.915070117300554e-08|n:3.166496753692627e-08|n:5.04078343549659e-08|n:0.01520044360584701|n:28.38566128414149|n:1.003013307112881|n:29.59409765142678|n:68.61646346676432|n:118.5903722863605|s:Hidden Size=[10], regularizer=0.05, learning_rate=0.1
n:8|s:model_2_8_14|n:0.999992187636361|n:0.9929823568752489|n:0.9999999792930766|n:0.9999999258332494|n:0.9999999621547906|n:2.949167273715139|n:2649.160279593511|n:0.008462013676762581|n:0.01518506277352571|n:0.01182353822514415|n:0.0006915313058469639|n:1.7173139706283|n:0.999937501090888|n:1.790423581687943|n:39.83695429918672|n:65.43334662141892|s:Hidden Size=[5], regularizer=0.2, learning_rate=0.1
n:9|s:model_30_8_1|n:0.9999798916619582|n:0.9920799592153482|n:0.9999912860188963|n:0.9999514375878201|n:0.9999912991856549|n:7.59089761078358|n:2989.815396206057|n:8.854970594868064|n:2.154577646637335|n:5.5047741207527|n:0.007582843109937573|n:2.755158364011692|n:1.00001930400452|n:2.872451159531733|n:93.94610030658913|n:153.671015725131|s:Hidden Size=[12], regularizer=0.2, learning_rate=0.1
n:10|s:model_8_7_0|n:0.9997581007167586|n:0.9942930704440628|n:0.9991719418988624|n:0.9999816159005624|n:0.9999890519567246|n:91.31697942361235|n:2154.365907366299|n:2.722241007490078|n:11.45177151635289|n:7.087006261921485|n:0.01626477638597697|n:9.555991807427022|n:1.005805582797793|n:9.962810162299261|n:40.97132651143687|n:71.44322213314189|s:Hidden Size=[6], regularizer=0.05, learning_rate=0.1
n:11|s:model_38_7_24|n:0.9996228785405012|n:0.9923654502610046|n:0.999951498895257|n:0.9999999999971473|n:0.9999609371212409|n:142.3633509607986|n:2882.042526470733|n:22.25018493272388|n:2.952292561531067e-07|n:11.12509261164826|n:0.03608018069230119|n:11.93161141509388|n:1.000274270152363|n:12.43956481487503|n:104.08323480199|n:173.5591568194774|s:Hidden Size=[14], regularizer=0.2, learning_rate=0.1
n:12|s:model_26_4_3|n:0.9994234533115562|n:0.9924212003350947|n:0.999913523046573|n:0.9999999999898611|n:0.9999201570513343|n:217.646374887526|n:2860.996873501741|n:40.49756265059114|n:1.867301762303528e-07|n:20.24916981847491|n:0.04177834619065011|n:14.75284294254928|n:1.000658910501079|n:15.38090200920857|n:79.23425677567501|n:134.083668894744|s:Hidden Size=[11], regularizer=0.2, learning_rate=0.1
n:13|s:model_35_7_20|n:0.9976680892119211|n:0.9907276899732727|n:0.9998407683893979|n:0.9999065008447229|n:0.9998789068593676|n:880.2963224998115|n:3500.297035089578|n:58.06455320422538|n:36.84860146441497|n:47.45597609574907|n:0.1152202243137621|n:29.66978804271799|n:1.001929857203927|n:30.93289234462473|n:92.43948283819037|n:157.039901556205|s:Hidden Size=[13], regularizer=0.05, learning_rate=0.01
n:14|s:model_18_8_24|n:0.9999120808381883|n:0.993368293836536|n:0.9998087303952133|n:0.9998677464268457|n:0.9998557449108948|n:33.18948358390484|n:2503.46907670764|n:56.39823384140618|n:150.1367359992585|n:103.2674849203323|n:0.01379699485226435|n:5.761031468748008|n:1.000162312298729|n:6.006290505351799|n:66.99553386755478|n:112.0939393876782|s:Hidden Size=[9], regularizer=0.2, learning_rate=0.1
n:15|s:model_118_8_22|n:0.9999976329375331|n:0.9930579544873804|n:0.9999987961061981|n:0.9999999910900658|n:0.9999999535882388|n:0.8935660812631268|n:2620.622181013896|n:0.05142131401225924|n:0.006023936904966831|n:0.02872262545861304|n:0.002968030022015975|n:0.9452862430307165|n:1.000056809499206|n:0.9855290354088171|n:50.22506997870457|n:80.69696560040958|s:Hidden Size=[3, 3], regularizer=0.2, learning_rate=0.1
n:16|s:model_21_9_0|n:0.9999197326609552|n:0.992652573302916|n:0.9998916382681973|n:0.9999221223374015|n:0.9999164517386357|n:30.30092048939316|n:2773.653578149225|n:98.80879856832325|n:9.176300647202879|n:53.99241118587088|n:0.03222094239256876|n:5.504627188955957|n:1.000275202305296|n:5.73897056453897|n:55.17764381751644|n:92.96279438843065|s:Hidden Size=[4, 3], regularizer=0.2, learning_rate=0.01
n:17|s:model_7_9_0|n:0.9998326480363232|n:0.9923832218103555|n:0.9999225595991194|n:0.9998541490986304|n:0.9998911204576901|n:63.1753662879765|n:2875.333766590793|n:13.79733842378482|n:39.84897306654602|n:26.82315574516542|n:0.04225721428789764|n:7.948293294033411|n:0.9991967105743514|n:8.286668594795767|n:29.70817109750642|n:52.86681177000224|s:Hidden Size=[2, 3], regularizer=0.05, learning_rate=0.01
n:18|s:model_29_6_0|n:0.9994984885356603|n:0.9938387948817379|n:0.9999486632114739|n:0.9998119424085081|n:0.9998516841850127|n:189.3205777882411|n:2325.85493214396|n:5.169935459503904|n:66.00821461365558|n:35.58907503657974|n:0.08809317910641791|n:13.7593814464256|n:1.000633488165482|n:14.3451468004461|n:75.51311648610645|n:127.9247769554391|s:Hidden Size=[2, 9], regularizer=0.2, learning_rate=0.01
"@

$lines = $data -split "`n"
$startRow = 9
for ($i = 0; $i -lt $lines.Length; $i++) {
    $line = $lines[$i]
    $fields = $line -split "\|"
    $rowNum = $startRow + $i
    for ($j = 0; $j -lt $fields.Length; $j++) {
        $field = $fields[$j]
        $tag = $field.Substring(0, 2)
        $val = $field.Substring(2)
        $addr = $cols[$j] + $rowNum
        if ($tag -eq "n:") {
            $ws.Range($addr).Value = [double]$val
        } else {
            $ws.Range($addr).Value = $val
        }
    }
}

Write-Output "Added $($lines.Length) rows"
